# Auto-generated script applying scheduled market-data refresh values
# to the Aegis_Profits workbook, per sheet/row/cell as captured in the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 123.166664
$ws.Range("I9").Value = 131.8
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 131.8
$ws.Range("L9").Value = 80
$ws.Range("M9").Value = 37.19999999999999
$ws.Range("N9").Value = -418

$ws.Range("H32").Value = 212.28572
$ws.Range("I32").Value = 296.33334
$ws.Range("J32").Value = 149.25
$ws.Range("K32").Value = 296.33334
$ws.Range("L32").Value = 149.25
$ws.Range("M32").Value = 29.66665999999998
$ws.Range("N32").Value = -801.25

$ws.Range("H70").Value = 1771.7222
$ws.Range("I70").Value = 1949.5714
$ws.Range("J70").Value = 1149.25
$ws.Range("K70").Value = 5848.7142
$ws.Range("L70").Value = 3447.75
$ws.Range("M70").Value = -5578.7142
$ws.Range("N70").Value = -3987.75

$ws.Range("H73").Value = 1771.7222
$ws.Range("I73").Value = 1949.5714
$ws.Range("J73").Value = 1149.25
$ws.Range("K73").Value = 5848.7142
$ws.Range("L73").Value = 3447.75
$ws.Range("M73").Value = -4912.7142
$ws.Range("N73").Value = -5319.75

$ws.Range("H98").Value = 1648.258
$ws.Range("I98").Value = 1446
$ws.Range("J98").Value = 2700
$ws.Range("K98").Value = 1446
$ws.Range("L98").Value = 2700
$ws.Range("M98").Value = 52
$ws.Range("N98").Value = -5696

$ws.Range("H114").Value = 66390
$ws.Range("J114").Value = 66390
$ws.Range("L114").Value = 66390
$ws.Range("N114").Value = -75068

$ws.Range("H122").Value = 1648.258
$ws.Range("I122").Value = 1446
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 4338
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -1888
$ws.Range("N122").Value = -13000

$ws.Range("H123").Value = 29818.092
$ws.Range("J123").Value = 29818.092
$ws.Range("L123").Value = 29818.092
$ws.Range("N123").Value = -39618.092

$ws.Range("H125").Value = 2768.1538
$ws.Range("I125").Value = 5881
$ws.Range("J125").Value = 2202.182
$ws.Range("K125").Value = 52929
$ws.Range("L125").Value = 19819.638
$ws.Range("M125").Value = -50469
$ws.Range("N125").Value = -24739.638

$ws.Range("H135").Value = 834.61536
$ws.Range("I135").Value = 708
$ws.Range("J135").Value = 4000
$ws.Range("K135").Value = 6372
$ws.Range("L135").Value = 36000
$ws.Range("M135").Value = -3837
$ws.Range("N135").Value = -41070

$ws.Range("H137").Value = 1619.326
$ws.Range("I137").Value = 1114.9487
$ws.Range("K137").Value = 3344.8461
$ws.Range("M137").Value = -794.8460999999998

$ws.Range("H138").Value = 5847.8506
$ws.Range("J138").Value = 9709.361000000001
$ws.Range("L138").Value = 29128.083
$ws.Range("N138").Value = -39408.083

$ws.Range("H141").Value = 1941.862
$ws.Range("I141").Value = 1661.8846
$ws.Range("K141").Value = 4985.6538
$ws.Range("M141").Value = 194.3462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 46737.184
$ws.Range("I45").Value = 67629.2
$ws.Range("J45").Value = 1968.5714
$ws.Range("K45").Value = 67629.2
$ws.Range("L45").Value = 1968.5714
$ws.Range("M45").Value = -67252.2
$ws.Range("N45").Value = -2722.5714

$ws.Range("H74").Value = 1986.8235
$ws.Range("I74").Value = 834.8
$ws.Range("J74").Value = 3632.5715
$ws.Range("K74").Value = 834.8
$ws.Range("L74").Value = 3632.5715
$ws.Range("M74").Value = 39.20000000000005
$ws.Range("N74").Value = -5380.5715

$ws.Range("H77").Value = 1986.8235
$ws.Range("I77").Value = 834.8
$ws.Range("J77").Value = 3632.5715
$ws.Range("K77").Value = 4174
$ws.Range("L77").Value = 18162.8575
$ws.Range("M77").Value = 194
$ws.Range("N77").Value = -26898.8575

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H82").Value = 41500
$ws.Range("J82").Value = 41500
$ws.Range("L82").Value = 41500
$ws.Range("N82").Value = -42222

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H85").Value = 41500
$ws.Range("J85").Value = 41500
$ws.Range("L85").Value = 41500
$ws.Range("N85").Value = -43996

$ws.Range("H86").Value = 49990
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 49990
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 49990
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -52362

$ws.Range("H89").Value = 49990
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 49990
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 149970
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -161826

$ws.Range("H122").Value = 1233.3334
$ws.Range("I122").Value = 1295.5714
$ws.Range("J122").Value = 1166.3077
$ws.Range("K122").Value = 3886.7142
$ws.Range("L122").Value = 3498.9231
$ws.Range("M122").Value = -1436.7142
$ws.Range("N122").Value = -8398.9231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 312.5
$ws.Range("I22").Value = 250
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 250
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -77
$ws.Range("N22").Value = -846

$ws.Range("H86").Value = 41942.07
$ws.Range("I86").Value = 63056.168
$ws.Range("J86").Value = 3936.7
$ws.Range("K86").Value = 63056.168
$ws.Range("L86").Value = 3936.7
$ws.Range("M86").Value = -61933.168
$ws.Range("N86").Value = -6182.7

$ws.Range("H89").Value = 41942.07
$ws.Range("I89").Value = 63056.168
$ws.Range("J89").Value = 3936.7
$ws.Range("K89").Value = 315280.84
$ws.Range("L89").Value = 19683.5
$ws.Range("M89").Value = -309664.84
$ws.Range("N89").Value = -30915.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30519.4
$ws.Range("I31").Value = 1179.5
$ws.Range("J31").Value = 42556.28
$ws.Range("K31").Value = 1179.5
$ws.Range("L31").Value = 42556.28
$ws.Range("M31").Value = -884.5
$ws.Range("N31").Value = -43146.28

$ws.Range("H34").Value = 30519.4
$ws.Range("I34").Value = 1179.5
$ws.Range("J34").Value = 42556.28
$ws.Range("K34").Value = 1179.5
$ws.Range("L34").Value = 42556.28
$ws.Range("M34").Value = -977.5
$ws.Range("N34").Value = -42960.28

$ws.Range("H92").Value = 22142.715
$ws.Range("J92").Value = 22142.715
$ws.Range("L92").Value = 22142.715
$ws.Range("N92").Value = -27134.715

$ws.Range("H132").Value = 20835408
$ws.Range("I132").Value = 17545814
$ws.Range("J132").Value = 33335866
$ws.Range("K132").Value = 52637442
$ws.Range("L132").Value = 100007598
$ws.Range("M132").Value = -52634912
$ws.Range("N132").Value = -100012658

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1638.2858
$ws.Range("I118").Value = 1078
$ws.Range("K118").Value = 3234
$ws.Range("M118").Value = -1991

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 91004550
$ws.Range("I80").Value = 143005980
$ws.Range("J80").Value = 2037.5
$ws.Range("K80").Value = 143005980
$ws.Range("L80").Value = 2037.5
$ws.Range("M80").Value = -143004982
$ws.Range("N80").Value = -4033.5

$ws.Range("H83").Value = 91004550
$ws.Range("I83").Value = 143005980
$ws.Range("J83").Value = 2037.5
$ws.Range("K83").Value = 715029900
$ws.Range("L83").Value = 10187.5
$ws.Range("M83").Value = -715024908
$ws.Range("N83").Value = -20171.5

$ws.Range("H122").Value = 2311
$ws.Range("I122").Value = 3750
$ws.Range("J122").Value = 1488.7142
$ws.Range("K122").Value = 11250
$ws.Range("L122").Value = 4466.142599999999
$ws.Range("M122").Value = -8800
$ws.Range("N122").Value = -9366.142599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1058.7142
$ws.Range("I22").Value = 1149.75
$ws.Range("J22").Value = 1037.2941
$ws.Range("K22").Value = 1149.75
$ws.Range("L22").Value = 1037.2941
$ws.Range("M22").Value = -854.75
$ws.Range("N22").Value = -1627.2941

$ws.Range("H27").Value = 1058.7142
$ws.Range("I27").Value = 1149.75
$ws.Range("J27").Value = 1037.2941
$ws.Range("K27").Value = 1149.75
$ws.Range("L27").Value = 1037.2941
$ws.Range("M27").Value = -1042.75
$ws.Range("N27").Value = -1251.2941

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H122").Value = 2494.1365
$ws.Range("I122").Value = 2186.5454
$ws.Range("J122").Value = 2801.7273
$ws.Range("K122").Value = 6559.6362
$ws.Range("L122").Value = 8405.1819
$ws.Range("M122").Value = -4109.6362
$ws.Range("N122").Value = -13305.1819

$ws.Range("H136").Value = 1697.258
$ws.Range("I136").Value = 1579.1364
$ws.Range("J136").Value = 1986
$ws.Range("K136").Value = 4737.4092
$ws.Range("L136").Value = 5958
$ws.Range("M136").Value = -2187.4092
$ws.Range("N136").Value = -11058

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 5671
$ws.Range("J41").Value = 6000
$ws.Range("L41").Value = 6000
$ws.Range("N41").Value = -6780

